$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "28.081.85"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -0.40%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.874.03"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -1.94%  "
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.22%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "313.48"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -0.42%  "
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +0.18%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5043"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  -0.53%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3834"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -2.17%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.08665"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -6.60%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.116"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -2.05%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "41.49"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -0.88%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "6.329"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -0.90%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "20.67"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -1.16%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "1.872.44"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -1.81%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "1.004"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +0.21%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "7.165"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -1.99%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001102"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -1.84%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "90.98"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -1.52%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06631"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +0.41%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "18.08"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +0.57%  "
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +0.21%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.100"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -2.00%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "28.119.56"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -0.47%  "
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +0.06%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.265"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -2.45%  "
$ws.Range("B26").NumberFormat = "@"
$ws.Range("B26").Value = "LidoDAOToken"
$ws.Range("C26").NumberFormat = "@"
$ws.Range("C26").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.564"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -1.00%  "
$ws.Range("B27").NumberFormat = "@"
$ws.Range("B27").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C27").NumberFormat = "@"
$ws.Range("C27").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.094.69"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -1.57%  "
$ws.Range("B28").NumberFormat = "@"
$ws.Range("B28").Value = "EthereumClassic"
$ws.Range("C28").NumberFormat = "@"
$ws.Range("C28").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "20.71"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -1.71%  "
$ws.Range("B29").NumberFormat = "@"
$ws.Range("B29").Value = "Monero"
$ws.Range("C29").NumberFormat = "@"
$ws.Range("C29").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "157.12"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -0.52%  "
$ws.Range("B30").NumberFormat = "@"
$ws.Range("B30").Value = "BitcoinCash"
$ws.Range("C30").NumberFormat = "@"
$ws.Range("C30").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "125.96"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -0.90%  "
$ws.Range("B31").NumberFormat = "@"
$ws.Range("B31").Value = "Stellar"
$ws.Range("C31").NumberFormat = "@"
$ws.Range("C31").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.1050"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -2.34%  "
$ws.Range("B32").NumberFormat = "@"
$ws.Range("B32").Value = "ImmutableX"
$ws.Range("C32").NumberFormat = "@"
$ws.Range("C32").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.060"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -3.50%  "
$ws.Range("B33").NumberFormat = "@"
$ws.Range("B33").Value = "Filecoin"
$ws.Range("C33").NumberFormat = "@"
$ws.Range("C33").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.593"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -0.78%  "
$ws.Range("B34").NumberFormat = "@"
$ws.Range("B34").Value = "HuobiToken"
$ws.Range("C34").NumberFormat = "@"
$ws.Range("C34").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.598"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -0.56%  "
$ws.Range("B35").NumberFormat = "@"
$ws.Range("B35").Value = "FraxShare"
$ws.Range("C35").NumberFormat = "@"
$ws.Range("C35").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "9.681"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -0.04%  "
$ws.Range("B36").NumberFormat = "@"
$ws.Range("B36").Value = "VeChain"
$ws.Range("C36").NumberFormat = "@"
$ws.Range("C36").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.02455"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +1.36%  "
$ws.Range("B37").NumberFormat = "@"
$ws.Range("B37").Value = "Hedera"
$ws.Range("C37").NumberFormat = "@"
$ws.Range("C37").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.06577"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -1.44%  "
$ws.Range("B38").NumberFormat = "@"
$ws.Range("B38").Value = "Algorand"
$ws.Range("C38").NumberFormat = "@"
$ws.Range("C38").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.2174"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -1.05%  "
$ws.Range("B39").NumberFormat = "@"
$ws.Range("B39").Value = "ARBITRUM"
$ws.Range("C39").NumberFormat = "@"
$ws.Range("C39").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.203"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -3.56%  "
$ws.Range("B40").NumberFormat = "@"
$ws.Range("B40").Value = "TrustWalletToken"
$ws.Range("C40").NumberFormat = "@"
$ws.Range("C40").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.246"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -3.49%  "
$ws.Range("B41").NumberFormat = "@"
$ws.Range("B41").Value = "Aptos"
$ws.Range("C41").NumberFormat = "@"
$ws.Range("C41").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "11.57"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +0.72%  "
$ws.Range("B42").NumberFormat = "@"
$ws.Range("B42").Value = "TheSandbox"
$ws.Range("C42").NumberFormat = "@"
$ws.Range("C42").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.6366"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -1.17%  "
$ws.Range("B43").NumberFormat = "@"
$ws.Range("B43").Value = "InternetComputer(DFINITY)"
$ws.Range("C43").NumberFormat = "@"
$ws.Range("C43").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "4.894"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -1.93%  "
$ws.Range("B44").NumberFormat = "@"
$ws.Range("B44").Value = "Frax"
$ws.Range("C44").NumberFormat = "@"
$ws.Range("C44").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.002"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +0.15%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "13.18"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -0.90%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.5982"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -0.81%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.281"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +0.09%  "
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -1.34%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.229"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +3.67%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.986"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -1.74%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "121.27"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -1.31%  "